# Edit script: resize/reposition two shapes on slide 1 (the rectangle and its
# matching picture) and split one italic text run on slide 2 into three runs
# (to wrap " * " in its own run), per the target OOXML diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 1: shrink + reposition the white rectangle (shape 4) and the picture
# behind it (shape 8) that together form the bottom-right chart tile.
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)

$rect = $slide1.Shapes.Item(4)
$rect.Left = 338.80517578125
$rect.Top = 153.7586212158203
$rect.Width = 294.14178466796875
$rect.Height = 177.77957153320312

$pic = $slide1.Shapes.Item(8)
$pic.Left = 338.8053894042969
$pic.Top = 153.75791931152344
$pic.Width = 294.180908203125
$pic.Height = 178.68310546875

# ---------------------------------------------------------------------------
# Slide 2: split the italic run "store’s mean monthly revenue * (state’s
# median rent/revenue ratio)" into three runs so " * " sits in its own run,
# keeping identical (italic, white, Oswald 12pt) formatting throughout.
# ---------------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$grp = $slide2.Shapes.Item(4)
$textShape = $grp.GroupItems.Item(2)
$tr = $textShape.TextFrame.TextRange

# Locate the italic run ("store's mean monthly revenue * ...") - it is the
# second distinct formatting run in the paragraph.
$italicRun = $tr.Runs(2, 1)

$italicRun.Text = "store’s mean monthly revenue"
$starRun = $italicRun.InsertAfter(" * ")
$null = $starRun.InsertAfter("(state’s median rent/revenue ratio)")

Write-Host "Edit applied."
